$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Igf2"
$ws.Range("C2").Value = "Insr"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 13.24090133333333
$ws.Range("H2").Value = 39.722704
$ws.Range("I2").Value = 0.1214410874295642
$ws.Range("J2").Value = 0.1214410874295642
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 18.444833
$ws.Range("N2").Value = 55.33449900000001
$ws.Range("O2").Value = 0.529296397589589
$ws.Range("P2").Value = 0.5292963975895891
$ws.Range("Q2").Value = 244.2262138628107
$ws.Range("R2").Value = 2198.035924765296
$ws.Range("S2").Value = 0.06427833009583067
$ws.Range("T2").Value = 0.06427833009583067

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Igf2"
$ws.Range("C3").Value = "Insr"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 13.24090133333333
$ws.Range("H3").Value = 39.722704
$ws.Range("I3").Value = 0.1214410874295642
$ws.Range("J3").Value = 0.1214410874295642
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 8.028768
$ws.Range("N3").Value = 24.086304
$ws.Range("O3").Value = 0.2303950368909585
$ws.Range("P3").Value = 0.2303950368909585
$ws.Range("Q3").Value = 106.308124916224
$ws.Range("R3").Value = 956.773124246016
$ws.Range("S3").Value = 0.02797942381841257
$ws.Range("T3").Value = 0.02797942381841257

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Igf2"
$ws.Range("C4").Value = "Insr"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 13.24090133333333
$ws.Range("H4").Value = 39.722704
$ws.Range("I4").Value = 0.1214410874295642
$ws.Range("J4").Value = 0.1214410874295642
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 8.374233
$ws.Range("N4").Value = 25.122699
$ws.Range("O4").Value = 0.2403085655194523
$ws.Range("P4").Value = 0.2403085655194524
$ws.Range("Q4").Value = 110.882392895344
$ws.Range("R4").Value = 997.9415360580961
$ws.Range("S4").Value = 0.02918333351532098
$ws.Range("T4").Value = 0.02918333351532098

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Igf2"
$ws.Range("C5").Value = "Insr"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 5.270503666666666
$ws.Range("H5").Value = 15.811511
$ws.Range("I5").Value = 0.04833928449947708
$ws.Range("J5").Value = 0.04833928449947708
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 18.444833
$ws.Range("N5").Value = 55.33449900000001
$ws.Range("O5").Value = 0.529296397589589
$ws.Range("P5").Value = 0.5292963975895891
$ws.Range("Q5").Value = 97.21355995755434
$ws.Range("R5").Value = 874.9220396179891
$ws.Range("S5").Value = 0.02558580914763148
$ws.Range("T5").Value = 0.02558580914763148

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Igf2"
$ws.Range("C6").Value = "Insr"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 5.270503666666666
$ws.Range("H6").Value = 15.811511
$ws.Range("I6").Value = 0.04833928449947708
$ws.Range("J6").Value = 0.04833928449947708
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 8.028768
$ws.Range("N6").Value = 24.086304
$ws.Range("O6").Value = 0.2303950368909585
$ws.Range("P6").Value = 0.2303950368909585
$ws.Range("Q6").Value = 42.315651182816
$ws.Range("R6").Value = 380.840860645344
$ws.Range("S6").Value = 0.01113713123553956
$ws.Range("T6").Value = 0.01113713123553956

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Igf2"
$ws.Range("C7").Value = "Insr"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 5.270503666666666
$ws.Range("H7").Value = 15.811511
$ws.Range("I7").Value = 0.04833928449947708
$ws.Range("J7").Value = 0.04833928449947708
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 8.374233
$ws.Range("N7").Value = 25.122699
$ws.Range("O7").Value = 0.2403085655194523
$ws.Range("P7").Value = 0.2403085655194524
$ws.Range("Q7").Value = 44.136425732021
$ws.Range("R7").Value = 397.227831588189
$ws.Range("S7").Value = 0.01161634411630603
$ws.Range("T7").Value = 0.01161634411630604

$ws.Range("A8").Value = "ECs"
$ws.Range("B8").Value = "Igf2"
$ws.Range("C8").Value = "Insr"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 90.52007366666668
$ws.Range("H8").Value = 271.560221
$ws.Range("I8").Value = 0.8302196280709587
$ws.Range("J8").Value = 0.8302196280709586
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 18.444833
$ws.Range("N8").Value = 55.33449900000001
$ws.Range("O8").Value = 0.529296397589589
$ws.Range("P8").Value = 0.5292963975895891
$ws.Range("Q8").Value = 1669.627641929365
$ws.Range("R8").Value = 15026.64877736428
$ws.Range("S8").Value = 0.4394322583461269
$ws.Range("T8").Value = 0.4394322583461269

$ws.Range("A9").Value = "ECs"
$ws.Range("B9").Value = "Igf2"
$ws.Range("C9").Value = "Insr"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 90.52007366666668
$ws.Range("H9").Value = 271.560221
$ws.Range("I9").Value = 0.8302196280709587
$ws.Range("J9").Value = 0.8302196280709586
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 8.028768
$ws.Range("N9").Value = 24.086304
$ws.Range("O9").Value = 0.2303950368909585
$ws.Range("P9").Value = 0.2303950368909585
$ws.Range("Q9").Value = 726.764670812576
$ws.Range("R9").Value = 6540.882037313184
$ws.Range("S9").Value = 0.1912784818370064
$ws.Range("T9").Value = 0.1912784818370064

$ws.Range("A10").Value = "ECs"
$ws.Range("B10").Value = "Igf2"
$ws.Range("C10").Value = "Insr"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 90.52007366666668
$ws.Range("H10").Value = 271.560221
$ws.Range("I10").Value = 0.8302196280709587
$ws.Range("J10").Value = 0.8302196280709586
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 8.374233
$ws.Range("N10").Value = 25.122699
$ws.Range("O10").Value = 0.2403085655194523
$ws.Range("P10").Value = 0.2403085655194524
$ws.Range("Q10").Value = 758.0361880618311
$ws.Range("R10").Value = 6822.325692556479
$ws.Range("S10").Value = 0.1995088878878253
$ws.Range("T10").Value = 0.1995088878878253

